# adding .rtf extension and styling
# Apply strike-through formatting to the phrase "modificarea fontului pe
# subsirurui" (and the single line break that immediately follows it),
# matching the rest of the already-struck-through text in that sentence.

$d = $word.ActiveDocument

$range = $d.Content
$found = $range.Find.Execute("modificarea fontului pe subsirurui", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $range.Font.StrikeThrough = $true

    # The phrase is immediately followed by a manual line break (<w:br/>)
    # that belongs to the same struck-through run in the source document.
    $lineBreak = $d.Range($range.End, $range.End + 1)
    $lineBreak.Font.StrikeThrough = $true
}
